# Weekly price-sheet update: two new daily records are inserted at the top
# of the data block (row 327), pushing all existing records down by two
# rows (old row 327 -> new row 329, ..., old row 425 -> new row 427).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 327/328; Excel shifts rows 327:425 down to
# 329:427 and extends the used range/dimension automatically.
$ws.Rows("327:328").Insert()

# New row 327: Poroto verde, Magnum, Primera — Región Metropolitana
$ws.Range("A327").Value = 9
$ws.Range("B327").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C327").Value = "Metropolitana"
$ws.Range("D327").Value = 44627
$ws.Range("E327").Value = 13
$ws.Range("F327").Value = 100112031
$ws.Range("G327").Value = "Poroto verde"
$ws.Range("H327").Value = "Magnum"
$ws.Range("I327").Value = "Primera"
$ws.Range("J327").Value = 79
$ws.Range("K327").Value = 25000
$ws.Range("L327").Value = 26000
$ws.Range("M327").Value = 25506
$ws.Range("N327").Value = "`$/saco 25 kilos"
$ws.Range("O327").Value = "Región Metropolitana"
$ws.Range("P327").Value = 1020
$ws.Range("Q327").Value = 25
$ws.Range("R327").Value = "Hortaliza"

# New row 328: Poroto verde, Sin especificar, Primera — Provincia del Elquí
$ws.Range("A328").Value = 9
$ws.Range("B328").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C328").Value = "Metropolitana"
$ws.Range("D328").Value = 44627
$ws.Range("E328").Value = 13
$ws.Range("F328").Value = 100112031
$ws.Range("G328").Value = "Poroto verde"
$ws.Range("H328").Value = "Sin especificar"
$ws.Range("I328").Value = "Primera"
$ws.Range("J328").Value = 34
$ws.Range("K328").Value = 34000
$ws.Range("L328").Value = 35000
$ws.Range("M328").Value = 34500
$ws.Range("N328").Value = "`$/malla 25 kilos"
$ws.Range("O328").Value = "Provincia del Elquí"
$ws.Range("P328").Value = 1380
$ws.Range("Q328").Value = 25
$ws.Range("R328").Value = "Hortaliza"
